$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '40.000.42'
$c.Style = 'Normal'

$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +1.47%  '
$c.Style = 'Normal'

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.200.78'
$c.Style = 'Normal'

$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +1.89%  '
$c.Style = 'Normal'

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.35%  '
$c.Style = 'Normal'

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '228.28'
$c.Style = 'Normal'

$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -0.31%  '
$c.Style = 'Normal'

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.628'
$c.Style = 'Normal'

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +1.11%  '
$c.Style = 'Normal'

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '63.51'
$c.Style = 'Normal'

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.54%  '
$c.Style = 'Normal'

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +0.13%  '
$c.Style = 'Normal'

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.396'
$c.Style = 'Normal'

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +0.07%  '
$c.Style = 'Normal'

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0856'
$c.Style = 'Normal'

$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -1.42%  '
$c.Style = 'Normal'

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +0.39%  '
$c.Style = 'Normal'

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '16.06'
$c.Style = 'Normal'

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +0.29%  '
$c.Style = 'Normal'

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '2.535.99'
$c.Style = 'Normal'

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +2.24%  '
$c.Style = 'Normal'

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '22.08'
$c.Style = 'Normal'

$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -0.42%  '
$c.Style = 'Normal'

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.818'
$c.Style = 'Normal'

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +0.31%  '
$c.Style = 'Normal'

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '5.58'
$c.Style = 'Normal'

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +0.19%  '
$c.Style = 'Normal'

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '2.212.33'
$c.Style = 'Normal'

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +2.07%  '
$c.Style = 'Normal'

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '40.003.54'
$c.Style = 'Normal'

$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +1.55%  '
$c.Style = 'Normal'

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.0₃0903'
$c.Style = 'Normal'

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +5.87%  '
$c.Style = 'Normal'

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '72.19'
$c.Style = 'Normal'

$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -0.16%  '
$c.Style = 'Normal'

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '6.07'
$c.Style = 'Normal'

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -1.09%  '
$c.Style = 'Normal'

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '231.58'
$c.Style = 'Normal'

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +1.25%  '
$c.Style = 'Normal'

$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.Style = 'Normal'

$c = $ws.Range('B24')
$c.NumberFormat = '@'
$c.Value = 'Toncoin'
$c.Style = 'Normal'

$c = $ws.Range('C24')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c.Style = 'Normal'

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.40'
$c.Style = 'Normal'

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +2.22%  '
$c.Style = 'Normal'

$c = $ws.Range('B25')
$c.NumberFormat = '@'
$c.Value = 'PancakeSwap'
$c.Style = 'Normal'

$c = $ws.Range('C25')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c.Style = 'Normal'

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.36'
$c.Style = 'Normal'

$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +0.17%  '
$c.Style = 'Normal'

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '9.67'
$c.Style = 'Normal'

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -0.88%  '
$c.Style = 'Normal'

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '171.67'
$c.Style = 'Normal'

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -0.13%  '
$c.Style = 'Normal'

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.140'
$c.Style = 'Normal'

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +1.59%  '
$c.Style = 'Normal'

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.46'
$c.Style = 'Normal'

$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +2.76%  '
$c.Style = 'Normal'

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '20.09'
$c.Style = 'Normal'

$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +2.08%  '
$c.Style = 'Normal'

$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +6.08%  '
$c.Style = 'Normal'

$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +1.07%  '
$c.Style = 'Normal'

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '4.59'
$c.Style = 'Normal'

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -1.66%  '
$c.Style = 'Normal'

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '4.73'
$c.Style = 'Normal'

$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  -1.47%  '
$c.Style = 'Normal'

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '7.05'
$c.Style = 'Normal'

$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -0.27%  '
$c.Style = 'Normal'

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.0622'
$c.Style = 'Normal'

$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +0.10%  '
$c.Style = 'Normal'

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '3.86'
$c.Style = 'Normal'

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +8.38%  '
$c.Style = 'Normal'

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.45'
$c.Style = 'Normal'

$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +1.08%  '
$c.Style = 'Normal'

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +0.29%  '
$c.Style = 'Normal'

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '4.96'
$c.Style = 'Normal'

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +17.53%  '
$c.Style = 'Normal'

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '103.58'
$c.Style = 'Normal'

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -0.24%  '
$c.Style = 'Normal'

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.0229'
$c.Style = 'Normal'

$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -0.82%  '
$c.Style = 'Normal'

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '17.75'
$c.Style = 'Normal'

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -1.90%  '
$c.Style = 'Normal'

$c = $ws.Range('B44')
$c.NumberFormat = '@'
$c.Value = 'Maker'
$c.Style = 'Normal'

$c = $ws.Range('C44')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c.Style = 'Normal'

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.518.96'
$c.Style = 'Normal'

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -0.96%  '
$c.Style = 'Normal'

$c = $ws.Range('B45')
$c.NumberFormat = '@'
$c.Value = 'TrustWalletToken'
$c.Style = 'Normal'

$c = $ws.Range('C45')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c.Style = 'Normal'

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '1.22'
$c.Style = 'Normal'

$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +2.93%  '
$c.Style = 'Normal'

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '8.25'
$c.Style = 'Normal'

$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +6.06%  '
$c.Style = 'Normal'

$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +0.14%  '
$c.Style = 'Normal'

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0925'
$c.Style = 'Normal'

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -0.78%  '
$c.Style = 'Normal'

$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -0.30%  '
$c.Style = 'Normal'

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.000194'
$c.Style = 'Normal'

$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +32.34%  '
$c.Style = 'Normal'

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.414.87'
$c.Style = 'Normal'

